# =====================================================================
# Edit: add "Player Info" and "ODI Batting Extra" sheets, rename the
# MATCH_CARD_LINK columns to MATCH_CODE in the existing sheets and
# replace the full scorecard URLs with the bare match code, and remove
# a few stray empty INNING_NUMBER cells.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Work out the existing sheets (before any insertions happen)
# ---------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------
# 2. Insert the new "Player Info" sheet before "ODI Batting"
# ---------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true

$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3447"
$playerInfo.Range("B2").Value = "Lendl Mark Platter Simmons"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------
# 3. Update the "ODI Batting" sheet
#    - rename MATCH_CARD_LINK -> MATCH_CODE
#    - replace the scorecard URL with the bare match code
#    - drop a handful of stray empty INNING_NUMBER cells
# ---------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLinkRange = $battingSheet.Range("D2:D69")
$battingLinkRange.NumberFormat = "@"
$battingLinkRange.Replace("http://www.howstat.com/cricket/Statistics/Matches/MatchScorecard_ODI.asp?MatchCode=", "")

$battingSheet.Range("B18").ClearContents()
$battingSheet.Range("B65").ClearContents()
$battingSheet.Range("B68").ClearContents()

# ---------------------------------------------------------------
# 4. Update the "ODI Bowling" sheet
#    - rename MATCH_CARD_LINK -> MATCH_CODE
#    - replace the scorecard URL with the bare match code
# ---------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLinkRange = $bowlingSheet.Range("B2:B13")
$bowlingLinkRange.NumberFormat = "@"
$bowlingLinkRange.Replace("http://www.howstat.com/cricket/Statistics/Matches/MatchScorecard_ODI.asp?MatchCode=", "")

# ---------------------------------------------------------------
# 5. Insert the new "ODI Batting Extra" sheet after "ODI Bowling"
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingExtra.Range("A1").Value = "MATCH_CODE"
$battingExtra.Range("B1").Value = "BATTING_POSITION"
$battingExtra.Range("C1").Value = "NUM_4"
$battingExtra.Range("D1").Value = "NUM_6"
$battingExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Range("F1").Value = "MAN_OF_MATCH"
$battingExtra.Range("A1:F1").Font.Bold = $true

# text formatting for the columns that must stay text even though some
# values look numeric (MATCH_CODE, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL,
# MAN_OF_MATCH). BATTING_POSITION (column B) is left as General so that
# it is stored as a genuine number.
$battingExtra.Range("A2:A21").NumberFormat = "@"
$battingExtra.Range("C2:F21").NumberFormat = "@"

$extraRows = @(
    @("3580", $null, $null, $null, $null,    "NO"),
    @("3581", $null, $null, $null, $null,    "NO"),
    @("3583", 5,     "2",   "0",   "4.94%",  "NO"),
    @("3593", $null, $null, $null, $null,    "NO"),
    @("3596", $null, $null, $null, $null,    "NO"),
    @("3597", 4,     "6",   "0",   "32.09%", "NO"),
    @("3598", $null, $null, $null, $null,    "NO"),
    @("3622", 5,     "6",   "0",   "24.16%", "NO"),
    @("3625", 5,     "4",   "2",   "44.03%", "NO"),
    @("3629", 4,     "3",   "0",   "5.76%",  "NO"),
    @("3655", 4,     "0",   "0",   $null,    "NO"),
    @("3657", 5,     "3",   "0",   "16.19%", "NO"),
    @("3661", $null, $null, $null, $null,    "NO"),
    @("3752", 6,     "9",   "5",   "33.55%", "NO"),
    @("3757", 6,     "4",   "2",   "16.13%", "NO"),
    @("3762", 6,     $null, $null, $null,    "NO"),
    @("3766", 6,     "0",   "0",   $null,    "NO"),
    @("3775", 6,     "0",   "0",   "4.95%",  "NO"),
    @("3788", $null, $null, $null, $null,    "NO"),
    @("3793", 3,     "1",   "1",   "4.80%",  "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $battingExtra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $battingExtra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $battingExtra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne $null) {
        $battingExtra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne $null) {
        $battingExtra.Cells.Item($r, 5).Value = $row[4]
    }
    $battingExtra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
